$wb = $excel.ActiveWorkbook

# Sheet "2025"
$ws = $wb.Worksheets.Item("2025")
$ws.Range("N2").Value = 7153.547888286129
$ws.Range("O2").Value = 6979.915717962018

# Sheet "2030"
$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 5707.815717280662
$ws.Range("I2").Value = 44492.05901988943
$ws.Range("L2").Value = 66334.06707325629
$ws.Range("M2").Value = 21991.42050229464
$ws.Range("O2").Value = 12078.10456148364

# Sheet "2035"
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 2927.360317916481
$ws.Range("B2").Value = 7940.887964949257
$ws.Range("E2").Value = 67179.99183625776
$ws.Range("I2").Value = 59530.75343380851
$ws.Range("L2").Value = 66334.06707325629
$ws.Range("M2").Value = 25547.11936466757
$ws.Range("N2").Value = 15114.04891232261
$ws.Range("O2").Value = 14759.64323227401

# Sheet "2040"
$ws = $wb.Worksheets.Item("2040")
$ws.Range("A2").Value = 2927.360317916481
$ws.Range("B2").Value = 7940.887964949257
$ws.Range("E2").Value = 67179.99183625776
$ws.Range("I2").Value = 59530.75343380851
$ws.Range("L2").Value = 66334.06707325629
$ws.Range("M2").Value = 25547.11936466757
$ws.Range("N2").Value = 15221.09667419122
$ws.Range("O2").Value = 14759.64323227401

# Sheet "2045"
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 6352.985609279765
$ws.Range("B2").Value = 7940.887964949257
$ws.Range("E2").Value = 67179.99183625776
$ws.Range("I2").Value = 59530.75343380851
$ws.Range("L2").Value = 66334.06707325629
$ws.Range("M2").Value = 25547.11936466757
$ws.Range("N2").Value = 15765.37192348151
$ws.Range("O2").Value = 17094.12375402333

# Sheet "2050"
$ws = $wb.Worksheets.Item("2050")
$ws.Range("A2").Value = 6352.985609279765
$ws.Range("B2").Value = 7940.887964949257
$ws.Range("E2").Value = 67179.99183625776
$ws.Range("I2").Value = 59530.75343380851
$ws.Range("L2").Value = 66334.06707325629
$ws.Range("M2").Value = 25547.11936466757
$ws.Range("N2").Value = 15765.37192348151
$ws.Range("O2").Value = 17094.12375402333
